$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay as literal text (avoid numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.867.90'
$ws.Range("D3").Value = '1.815.22'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '308.86'
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4661'
$ws.Range("E7").Value = '  +1.68%  '
$ws.Range("D8").Value = '0.3684'
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("E9").Value = '  +1.74%  '
$ws.Range("D10").Value = '0.8701'
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '1.805.21'
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").Value = '5.360'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").Value = '0.07065'
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").Value = '6.498'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '91.49'
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = '0.000008689'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '14.73'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("D21").Value = '26.903.79'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").Value = '5.336'
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").Value = '2.108.25'
$ws.Range("E24").Value = '  +3.74%  '
$ws.Range("D25").Value = '1.902'
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("D26").Value = '150.18'
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("D27").Value = '2.177'
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").Value = '18.31'
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("D29").Value = '5.315'
$ws.Range("E29").Value = '  +2.03%  '
$ws.Range("D30").Value = '115.71'
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").Value = '0.08934'
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D32").Value = '0.7667'
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("D33").Value = '1.163'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").Value = '4.506'
$ws.Range("E34").Value = '  +1.78%  '
$ws.Range("D35").Value = '2.900'
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '1.085'
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").Value = '0.05284'
$ws.Range("E39").Value = '  +1.49%  '
$ws.Range("D40").Value = '2.941'
$ws.Range("E40").Value = '  +1.48%  '
$ws.Range("D41").Value = '7.257'
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").Value = '0.5311'
$ws.Range("E42").Value = '  +1.83%  '
$ws.Range("D43").Value = '2.343'
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").Value = '0.1662'
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("D45").Value = '8.417'
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").Value = '0.4931'
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.669'
$ws.Range("E49").Value = '  +1.68%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '103.76'
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").Value = '0.06285'
$ws.Range("E51").Value = '  -0.32%  '
